# Insert a new weekly price-report row above row 34 ("Fruta / hortaliza, semanal"):
# all the existing rows from 34 down get pushed one row lower (34->35, 35->36, ...,
# 145->146) and the freshly-opened row 34 is filled in with the new week's data,
# using the same Mercado/Región/Categoría/Unidad/Origen/Clasificación boilerplate
# that every other row in this sheet already shares.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 34..145 down to 35..146, opening up a blank row 34.
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new week's record.
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44707
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112044
$ws.Range("G34").Value = "Perejil"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1750
$ws.Range("N34").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O34").Value = "Provincia del Elquí"
$ws.Range("P34").Value = 1167
$ws.Range("Q34").Value = 1.5
$ws.Range("R34").Value = "Hortaliza"
